# Apply the changes described by the commit:
# - Add a new "Mujeres" worksheet after "Mujeres<1" with the infant-girls data table
# - Update selections/active sheet to reflect the edited workbook state

$wb = $excel.ActiveWorkbook

$sheetHombresLt1 = $wb.Worksheets.Item(1)   # "Hombres<1"
$sheetHombres    = $wb.Worksheets.Item(2)   # "Hombres"
$sheetMujeresLt1 = $wb.Worksheets.Item(3)   # "Mujeres<1"

# ---------------------------------------------------------------------------
# 1. Create the new "Mujeres" sheet, placed right after "Mujeres<1"
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $sheetMujeresLt1)
$newSheet.Name = "Mujeres"

# Header row
$newSheet.Range("A1").Value = "Edad (años)"
$newSheet.Range("B1").Value = "Peso (Kg)"
$newSheet.Range("C1").Value = "Talla (cm)"

# Data rows
$newSheet.Range("A2").Value = 7
$newSheet.Range("B2").Value = 10
$newSheet.Range("A3").Value = 8
$newSheet.Range("C3").Value = 155

# Copy the cell formatting (styles) from the matching cells on the sibling
# sheets so the new sheet keeps the same visual style as the rest of the
# workbook.
$sheetMujeresLt1.Range("A1:B2").Copy()
$newSheet.Range("A1:B2").PasteSpecial(-4122)

$sheetMujeresLt1.Range("B1").Copy()
$newSheet.Range("C1").PasteSpecial(-4122)

$sheetMujeresLt1.Range("B2").Copy()
$newSheet.Range("C3").PasteSpecial(-4122)

$sheetHombresLt1.Range("A3").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Update selections on each sheet
# ---------------------------------------------------------------------------

# "Hombres" selection moves to C1:C1048576 (whole column C selected)
$sheetHombres.Activate()
$sheetHombres.Range("C1:C1048576").Select()

# "Mujeres<1" selection moves to G23
$sheetMujeresLt1.Activate()
$sheetMujeresLt1.Range("G23").Select()

# New "Mujeres" sheet selection is C10, and it becomes the active tab
$newSheet.Activate()
$newSheet.Range("C10").Select()
